# This workbook ("Example_Marks.xlsx") is being turned into a generic
# example/template used by a new GUI feedback-form generator. The real
# student rows and their long personalised feedback comments are trimmed
# down to a couple of placeholder rows with short, generic example text,
# and the leftover (now unused) rows are blanked out completely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2 and 3 ("Cater, Jack" / "Crick, Jo") keep their Name/Mark values but
# their long, student-specific feedback text is replaced with a short
# generic placeholder.
$ws.Range("C2").Value = "Example Feedback"
$ws.Range("C3").Value = "Example Feedback"

# Rows 4-8 ("Owens, Mike", "Shirley, Poop", "Ure, Alper", "Iqbal, Jackie",
# "James, Dean") are no longer needed as example data, so wipe their
# content and formatting entirely, leaving just blank rows (row heights
# are preserved automatically since we only clear the cell range).
$ws.Range("A4:D8").Clear()

# Move/save the active selection to B5, matching where the author left the
# cursor when they saved the trimmed-down example workbook.
$ws.Range("B5").Select()
